# This workbook holds one weekly price observation (a "Primera"/"Segunda"
# row pair) per data-refresh, appended at the top of the historical block
# that starts at row 121 (rows 1-120 are the header + the untouched
# earlier weeks).
#
# The new refresh inserts a brand-new pair of rows for the latest week
# (date serial 44512) right above the existing block, which pushes every
# row from the old 121 down to 123, old 122 down to 124, ..., old 222
# down to 224. The new top rows (121/122) repeat the same
# category/quality/price layout that used to sit in that slot, just with
# the refreshed date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: push rows 121:222 down to 123:224 by inserting two blank
#    rows at 121. Excel's Insert shifts everything below (and grows the
#    used range / <dimension> automatically).
$ws.Rows("121:122").Insert()

# 2) The rows that used to be 121:122 are now sitting at 123:124 (shifted
#    down by the insert). Clone that pair back into the freshly-opened
#    121:122 slot so the new entry starts from the same template.
$ws.Range("A123:R124").Copy()
$ws.Range("A121").PasteSpecial()
$excel.CutCopyMode = $false

# 3) Stamp the new pair with the latest date (2021-11-12 == serial 44512),
#    while every other column keeps the template values just pasted in.
$ws.Range("D121").Value2 = 44512
$ws.Range("D122").Value2 = 44512
